$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = $origStyle
}

Set-TextValue "D2" "26.717.47"
$ws.Range("E2").Value = "  +1.42%  "
Set-TextValue "D3" "1.634.53"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  +0.18%  "
Set-TextValue "D5" "213.77"
$ws.Range("E5").Value = "  +0.75%  "
Set-TextValue "D6" "0.500"
$ws.Range("E6").Value = "  +3.40%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("E9").Value = "  +1.02%  "
Set-TextValue "D10" "19.27"
$ws.Range("E11").Value = "  +3.52%  "
$ws.Range("E12").Value = "  +0.69%  "
Set-TextValue "D13" "1.626.96"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("E15").Value = "  +0.86%  "
Set-TextValue "D16" "26.705.09"
$ws.Range("E16").Value = "  +1.30%  "
Set-TextValue "D17" "63.63"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("E18").Value = "  +2.20%  "
Set-TextValue "D19" "218.99"
$ws.Range("E19").Value = "  +8.37%  "
$ws.Range("E20").Value = "  +0.20%  "
Set-TextValue "D21" "4.32"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("E22").Value = "  +0.77%  "
Set-TextValue "D23" "6.16"
$ws.Range("E23").Value = "  +1.72%  "
Set-TextValue "D24" "1.97"
$ws.Range("E24").Value = "  +5.38%  "
Set-TextValue "D25" "147.31"
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("E28").Value = "  +4.11%  "
Set-TextValue "D29" "15.57"
$ws.Range("E29").Value = "  +2.51%  "
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  +3.65%  "
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("E34").Value = "  +0.77%  "
Set-TextValue "D35" "1.229.07"
$ws.Range("E35").Value = "  +5.89%  "
$ws.Range("E36").Value = "  +0.14%  "
Set-TextValue "D37" "0.0173"
$ws.Range("E37").Value = "  +5.75%  "
Set-TextValue "D38" "0.808"
$ws.Range("E38").Value = "  +0.78%  "
Set-TextValue "D40" "0.502"
$ws.Range("E40").Value = "  +1.32%  "
Set-TextValue "D41" "2.29"
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("E43").Value = "  -0.59%  "
Set-TextValue "D44" "1.769.82"
$ws.Range("E44").Value = "  +0.54%  "
Set-TextValue "D45" "92.80"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("E46").Value = "  +3.04%  "
Set-TextValue "D47" "55.45"
$ws.Range("E47").Value = "  +2.81%  "
$ws.Range("E48").Value = "  -0.42%  "
Set-TextValue "D49" "0.0513"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("E50").Value = "  +4.74%  "
$ws.Range("E51").Value = "  -0.27%  "
